$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new analyte below the existing list (A1:A15 -> A1:A16)
$ws.Range("A16").Value = "test"

# Match the formatting used by the other analyte rows (A2:A15) by copying
# the format of the row directly above down onto the new row.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Rows("16").RowHeight = 15.75

$excel.CutCopyMode = $false

$ws.Range("A16").Select()
